# Parliament_Members.xlsx update
#  - Update 6 "Topic opinion" sentences on the Content sheet (G7, G8, G9, G10, G11, G13)
#    with new "negative" phrasing (adds 6 new shared strings).
#  - Mark two cells (I5, I6) on the Content sheet with the existing date number format
#    (empty placeholder cells, same style as the Publication Date column).
#  - Add a new "Votes" worksheet (after "Content") listing each member's three recorded
#    votes (by date) on the "Ban of Thermal Vehicles" topic.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Content sheet: replace a handful of "Published Content" opinions with
#    new negative-themed statements.
# ---------------------------------------------------------------------------
$wsContent = $wb.Worksheets.Item("Content")

$wsContent.Range("G7").Value2  = "New measures are not needed to support industries affected by the transition."
$wsContent.Range("G8").Value2  = "I believe this policy is not necessary to combat climate change."
$wsContent.Range("G9").Value2  = "The current measures are enough to support industries affected by the transition."
$wsContent.Range("G10").Value2 = "Public opinion is stable so nothing needs to change."
$wsContent.Range("G11").Value2 = "Climate change is a hoax."
$wsContent.Range("G13").Value2 = "I don't believe this policy is necessary to combat climate change."

# Empty cells in column I (rows 5-6) carrying the existing date-number format.
$wsContent.Range("I5:I6").NumberFormat = "yyyy\-mm\-dd"

# ---------------------------------------------------------------------------
# 2) Add the new "Votes" worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
$wsMembers = $wb.Worksheets.Item("Parliament_Members")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVotes = $wb.Worksheets.Add($null, $lastSheet)
$wsVotes.Name = "Votes"

$wsVotes.Range("A1").Value2 = "ID"
$wsVotes.Range("B1").Value2 = "Last Name"
$wsVotes.Range("C1").Value2 = "First Name"
$wsVotes.Range("D1").Value2 = 45323
$wsVotes.Range("E1").Value2 = 45439
$wsVotes.Range("F1").Value2 = 45610

# Copy the bold/bordered/centered header formatting already used for the
# other sheets' header rows onto A1:C1.
$wsContent.Range("A1:C1").Copy()
$wsVotes.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header date cells: bold + the existing date number format (creates the new,
# bold variant of that number format).
$wsVotes.Range("D1:F1").NumberFormat = "yyyy\-mm\-dd"
$wsVotes.Range("D1:F1").Font.Bold = $true

# Member votes, mirroring Parliament_Members' Last/First Name order.
$votes = @(
    @(1,  "Richard",    "Katherine",   1, 1, 1),
    @(2,  "Sullivan",   "Kristen",     0, 0, 0),
    @(3,  "Wagner",     "Juan",        0, 1, 1),
    @(4,  "Mcgrath",    "Lisa",        0, 0, $null),
    @(5,  "Nolan",      "Gabriella",   0, 0, $null),
    @(6,  "Stewart",    "Danielle",    0, 0, $null),
    @(7,  "Green",      "Linda",       0, 0, $null),
    @(8,  "Page",       "Kristen",     0, 0, $null),
    @(9,  "Sanchez",    "Mark",        0, 0, $null),
    @(10, "Pratt",      "Douglas",     0, 0, $null),
    @(11, "Golden",     "Whitney",     0, 0, $null),
    @(12, "Mann",       "Amy",         0, 0, $null),
    @(13, "Roberts",    "Christopher", 0, 0, $null),
    @(14, "Roberts",    "Jaime",       0, 0, $null),
    @(15, "Alvarez",    "Christopher", 0, 0, $null),
    @(16, "Villanueva", "Sandra",      0, 0, $null),
    @(17, "Hill",       "Lisa",        0, 0, $null),
    @(18, "Hill",       "Travis",      0, 0, $null),
    @(19, "Nguyen",     "Darrell",     0, 0, $null),
    @(20, "Smith",      "Joseph",      0, 0, $null)
)

$r = 2
foreach ($row in $votes) {
    $wsVotes.Cells.Item($r, 1).Value2 = $row[0]
    $wsVotes.Cells.Item($r, 2).Value2 = $row[1]
    $wsVotes.Cells.Item($r, 3).Value2 = $row[2]
    $wsVotes.Cells.Item($r, 4).Value2 = $row[3]
    $wsVotes.Cells.Item($r, 5).Value2 = $row[4]
    if ($row[5] -ne $null) {
        $wsVotes.Cells.Item($r, 6).Value2 = $row[5]
    }
    $r++
}

# ---------------------------------------------------------------------------
# 3) Restore/update the on-screen selections so they match the edited file.
#    "Content" is left as the active/front-most tab, so it is selected last.
# ---------------------------------------------------------------------------
$wsMembers.Range("A1:C21").Select()

$wsVotes.Activate()
$wsVotes.Range("F5").Select()

$wsContent.Activate()
$wsContent.Range("B2:F2").Select()

Write-Output "done"
